$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.944.06'
$ws.Range('E2').Value = '  -1.43%  '
$ws.Range('D3').Value = '2.282.30'
$ws.Range('E3').Value = '  -3.08%  '
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '311.06'
$ws.Range('E5').Value = '  -4.09%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '105.91'
$ws.Range('E6').Value = '  +3.64%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.626'
$ws.Range('E7').Value = '  -1.90%  '
$ws.Range('E8').Value = '  +0.12%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.605'
$ws.Range('E9').Value = '  -2.93%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '40.26'
$ws.Range('E10').Value = '  +0.13%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0908'
$ws.Range('E11').Value = '  -1.53%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '8.26'
$ws.Range('E12').Value = '  -2.13%  '
$ws.Range('E13').Value = '  +0.00%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.965'
$ws.Range('E14').Value = '  -3.30%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.45'
$ws.Range('E15').Value = '  -4.17%  '
$ws.Range('D16').Value = '2.622.21'
$ws.Range('E16').Value = '  -3.15%  '
$ws.Range('D17').Value = '2.285.26'
$ws.Range('E17').Value = '  -2.48%  '
$ws.Range('D18').Value = '41.871.91'
$ws.Range('E18').Value = '  -1.82%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.54'
$ws.Range('E19').Value = '  -4.67%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0000105'
$ws.Range('E20').Value = '  -1.81%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '73.36'
$ws.Range('E21').Value = '  -4.20%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.43'
$ws.Range('E22').Value = '  -6.87%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '255.70'
$ws.Range('E23').Value = '  -2.93%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.31'
$ws.Range('E24').Value = '  +0.20%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.30'
$ws.Range('E25').Value = '  -7.13%  '
$ws.Range('E26').Value = '  +0.48%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.92'
$ws.Range('E27').Value = '  -4.71%  '
$ws.Range('E28').Value = '  +3.33%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '22.58'
$ws.Range('E29').Value = '  -0.53%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '166.06'
$ws.Range('E30').Value = '  -5.30%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '35.57'
$ws.Range('E31').Value = '  +0.60%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0887'
$ws.Range('E32').Value = '  -1.22%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.92'
$ws.Range('E33').Value = '  -5.91%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.75'
$ws.Range('E34').Value = '  -4.68%  '
$ws.Range('B35').Value = 'Kaspa'
$ws.Range('C35').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.118'
$ws.Range('E35').Value = '  +7.80%  '
$ws.Range('B36').Value = 'Stellar'
$ws.Range('C36').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.129'
$ws.Range('E36').Value = '  -2.42%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.55'
$ws.Range('E37').Value = '  -0.02%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0352'
$ws.Range('E38').Value = '  -1.47%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.78'
$ws.Range('E39').Value = '  -0.62%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.62'
$ws.Range('E40').Value = '  -4.25%  '
$ws.Range('B41').Value = 'BitcoinSV'
$ws.Range('C41').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '98.57'
$ws.Range('E41').Value = '  +6.79%  '
$ws.Range('B42').Value = 'MultiversX'
$ws.Range('C42').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '71.46'
$ws.Range('E42').Value = '  +2.24%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.46'
$ws.Range('E43').Value = '  -2.57%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.227'
$ws.Range('E44').Value = '  -4.53%  '
$ws.Range('E45').Value = '  +0.13%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '12.25'
$ws.Range('E46').Value = '  +2.93%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '112.03'
$ws.Range('E47').Value = '  -7.78%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.03'
$ws.Range('E48').Value = '  -1.43%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '5.28'
$ws.Range('E49').Value = '  -4.39%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '75.10'
$ws.Range('E50').Value = '  +6.38%  '
$ws.Range('D51').Value = '1.558.87'
$ws.Range('E51').Value = '  +0.92%  '
